$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) and Volume(1h) (E) columns contain numeric-looking text values
# in the source workbook (stored as text, not numbers). Force the Text number
# format on each cell before assigning so Excel keeps them as text strings
# instead of converting to floats/percentages.

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2.12%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "38.88"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "8.72%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.123"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.75%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08185"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.82%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.003"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "7.23%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.197"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2.01%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "7.922"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2.03%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9326"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.38%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1406"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "5.08%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1957"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "3.42%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09147"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.26%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03469"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.12%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09850"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.13%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001411"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.96%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005825"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-4.87%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.568"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-4.51%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "3.89%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3450"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.18%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1338"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.13%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.814"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-6.60%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2470"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "5.17%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04470"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.07%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001241"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.29%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001303"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.30%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02115"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "8.72%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05186"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-4.18%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007480"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.72%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01002"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-1.08%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1369"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "1.34%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002135"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.17%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009770"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-3.97%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006327"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.87%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000752"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.40%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-0.63%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-3.11%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002105"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.40%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002005"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.40%"

# Coin name (B) and Link (C) columns are plain text; assigning directly keeps them as strings.
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
